$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("O3").Value = 1.21
$ws.Range("R3").Value = 1.56
$ws.Range("AA4").Value = 1000
$ws.Range("AB4").Value = 1000
$ws.Range("AC4").Value = 1000
$ws.Range("AD4").Value = 1000
$ws.Range("AE4").Value = 1000
$ws.Range("AF4").Value = 1000
$ws.Range("AG4").Value = 1000
$ws.Range("AH4").Value = 1000
$ws.Range("AI4").Value = 1000
$ws.Range("AJ4").Value = 1000
$ws.Range("AK4").Value = 1000
$ws.Range("AL4").Value = 1000
$ws.Range("AM4").Value = 1000
$ws.Range("AN4").Value = 1000
$ws.Range("AO4").Value = 1000
$ws.Range("F4").Value = 1.89
$ws.Range("H4").Value = 2.84
$ws.Range("K4").Value = 1000
$ws.Range("L4").Value = 1.01
$ws.Range("M4").Value = 1.01
$ws.Range("R4").Value = 1.24
$ws.Range("S4").Value = 3.1
$ws.Range("T4").Value = 1.01
$ws.Range("U4").Value = 1.01
$ws.Range("V4").Value = 1.36
$ws.Range("W4").Value = 1.56
$ws.Range("X4").Value = 1000
$ws.Range("Y4").Value = 1000
$ws.Range("Z4").Value = 1000
$ws.Range("G5").Value = 3.75
$ws.Range("J5").Value = 3.2
$ws.Range("K5").Value = 3.25
$ws.Range("F6").Value = 1.78
$ws.Range("I6").Value = 5.3
$ws.Range("G7").Value = 6.2
$ws.Range("H7").Value = 1.61
$ws.Range("T7").Value = 1.73
$ws.Range("F8").Value = 2.98
$ws.Range("G8").Value = 3.2
$ws.Range("H8").Value = 2.54
$ws.Range("I8").Value = 2.6
$ws.Range("K12").Value = 3.65
$ws.Range("N13").Value = 3.25
$ws.Range("O13").Value = 1.42
$ws.Range("R13").Value = 1.28
$ws.Range("T13").Value = 2.02
$ws.Range("X13").Value = 11.5
$ws.Range("N19").Value = 3.55
$ws.Range("F20").Value = 1.78
$ws.Range("G20").Value = 2.26
$ws.Range("H20").Value = 3.25
$ws.Range("I20").Value = 7
$ws.Range("J20").Value = 3.4
$ws.Range("K20").Value = 7.8
$ws.Range("P20").Value = 2.14
$ws.Range("AO23").Value = 5.6
$ws.Range("F23").Value = 8.199999999999999
$ws.Range("H23").Value = 1.39
$ws.Range("I23").Value = 1.44
$ws.Range("K23").Value = 5.8
$ws.Range("R23").Value = 1.59
$ws.Range("F24").Value = 1.81
$ws.Range("G24").Value = 1.85
$ws.Range("H24").Value = 4.9
$ws.Range("P24").Value = 1.9
$ws.Range("Q24").Value = 1.97
$ws.Range("I25").Value = 9.199999999999999
$ws.Range("J25").Value = 5.2
$ws.Range("K25").Value = 5.3
$ws.Range("F26").Value = 1.85
$ws.Range("K26").Value = 4.1
$ws.Range("AC28").Value = 17
$ws.Range("AF28").Value = 7.8
$ws.Range("AJ28").Value = 9
$ws.Range("AL28").Value = 38
$ws.Range("AM28").Value = 180
$ws.Range("H28").Value = 13
$ws.Range("I28").Value = 13.5
$ws.Range("Q28").Value = 1.55
$ws.Range("N29").Value = 5
$ws.Range("P29").Value = 2.32
$ws.Range("X29").Value = 19
$ws.Range("N30").Value = 3.15
$ws.Range("Q30").Value = 2.34
$ws.Range("T30").Value = 2.06
$ws.Range("U30").Value = 1.9
$ws.Range("H34").Value = 2.56
$ws.Range("I34").Value = 2.6
$ws.Range("K34").Value = 3.25
$ws.Range("AG35").Value = 10.5
$ws.Range("AH35").Value = 16.5
$ws.Range("AK35").Value = 16
$ws.Range("I35").Value = 4.3
$ws.Range("S35").Value = 2.88
$ws.Range("F36").Value = 1.52
$ws.Range("P36").Value = 2.18
$ws.Range("AC39").Value = 10.5
$ws.Range("J39").Value = 4.6
$ws.Range("K39").Value = 4.8
$ws.Range("Y39").Value = 7.2
